$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.420.72'
$ws.Range("E2").Value = '  +0.06%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.668.12'
$ws.Range("E3").Value = '  -0.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '642.39'
$ws.Range("E5").Value = '  -5.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.59'
$ws.Range("E6").Value = '  -0.50%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.496'
$ws.Range("E8").Value = '  +0.32%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.144'
$ws.Range("E9").Value = '  -1.20%  '

# Row 10
$ws.Range("E10").Value = '  -0.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.438'
$ws.Range("E11").Value = '  +0.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000229'
$ws.Range("E12").Value = '  -1.24%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.292.04'
$ws.Range("E13").Value = '  -0.44%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.26'
$ws.Range("E14").Value = '  -0.65%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.671.60'
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.471.83'
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("E17").Value = '  +1.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.96'
$ws.Range("E18").Value = '  -0.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.41'
$ws.Range("E19").Value = '  -0.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.52'
$ws.Range("E20").Value = '  -0.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.69'
$ws.Range("E21").Value = '  -3.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.641'
$ws.Range("E22").Value = '  -1.67%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.40'
$ws.Range("E23").Value = '  -0.61%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.817.72'

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +0.49%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.77'
$ws.Range("E27").Value = '  -1.69%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.87'
$ws.Range("E28").Value = '  -3.08%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.59'
$ws.Range("E29").Value = '  -3.31%  '

# Row 30
$ws.Range("E30").Value = '  -6.20%  '

# Row 31
$ws.Range("E31").Value = '  +0.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.98'
$ws.Range("E32").Value = '  -0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.52'
$ws.Range("E33").Value = '  -1.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.41'
$ws.Range("E34").Value = '  -3.39%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.663.44'

# Row 36
$ws.Range("E36").Value = '  +1.85%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.31'
$ws.Range("E37").Value = '  +0.33%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '179.49'
$ws.Range("E39").Value = '  +4.47%  '

# Row 40
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("E40").Value = '  -5.48%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0891'
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.17'
$ws.Range("E43").Value = '  -3.72%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.927'
$ws.Range("E44").Value = '  -1.77%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.69'
$ws.Range("E45").Value = '  -2.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.69'
$ws.Range("E46").Value = '  -0.88%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.84'
$ws.Range("E47").Value = '  -5.49%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.25'
$ws.Range("E48").Value = '  -3.04%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.78'
$ws.Range("E49").Value = '  -0.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000262'
$ws.Range("E50").Value = '  -5.27%  '

# Row 51
$ws.Range("E51").Value = '  -5.76%  '
